$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 339, shifting the existing rows 339-366 down to 341-368
$ws.Rows("339:340").Insert()

# Fill in the new row 339 (Zanahoria, Primera, Araucania)
$ws.Cells.Item(339, 1).Value = 11
$ws.Cells.Item(339, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(339, 3).Value = "Bíobío"
$ws.Cells.Item(339, 4).Value = 45013
$ws.Cells.Item(339, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(339, 5).Value = 8
$ws.Cells.Item(339, 6).Value = 100114013
$ws.Cells.Item(339, 7).Value = "Zanahoria"
$ws.Cells.Item(339, 8).Value = "Sin especificar"
$ws.Cells.Item(339, 9).Value = "Primera"
$ws.Cells.Item(339, 10).Value = 200
$ws.Cells.Item(339, 11).Value = 3500
$ws.Cells.Item(339, 12).Value = 4000
$ws.Cells.Item(339, 13).Value = 3750
$ws.Cells.Item(339, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(339, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(339, 16).Value = 188
$ws.Cells.Item(339, 17).Value = 20
$ws.Cells.Item(339, 18).Value = "Hortaliza"

# Fill in the new row 340 (Zanahoria, Primera, Nuble)
$ws.Cells.Item(340, 1).Value = 11
$ws.Cells.Item(340, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(340, 3).Value = "Bíobío"
$ws.Cells.Item(340, 4).Value = 45013
$ws.Cells.Item(340, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(340, 5).Value = 8
$ws.Cells.Item(340, 6).Value = 100114013
$ws.Cells.Item(340, 7).Value = "Zanahoria"
$ws.Cells.Item(340, 8).Value = "Sin especificar"
$ws.Cells.Item(340, 9).Value = "Primera"
$ws.Cells.Item(340, 10).Value = 270
$ws.Cells.Item(340, 11).Value = 5000
$ws.Cells.Item(340, 12).Value = 6000
$ws.Cells.Item(340, 13).Value = 5556
$ws.Cells.Item(340, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(340, 15).Value = "Región de Ñuble"
$ws.Cells.Item(340, 16).Value = 278
$ws.Cells.Item(340, 17).Value = 20
$ws.Cells.Item(340, 18).Value = "Hortaliza"

Write-Host "Done inserting rows 339-340"
